$wb = $excel.ActiveWorkbook

# --- Sheet "Assign 2": append the "8th Attempt" block (rows 38-41) ---
$ws2 = $wb.Worksheets.Item("Assign 2")

# Copy formatting down from the matching "7th Attempt" rows (33-36) so the
# new block gets the same cell styles (s="1" for A/B, s="2" for C:F).
$ws2.Range("A33:G33").Copy()
$ws2.Range("A38").PasteSpecial(-4122)

$ws2.Range("A34:F34").Copy()
$ws2.Range("A39").PasteSpecial(-4122)

$ws2.Range("A35:F35").Copy()
$ws2.Range("A40").PasteSpecial(-4122)

$ws2.Range("A36:F36").Copy()
$ws2.Range("A41").PasteSpecial(-4122)

# Row 38: Original / Linear Regression
$ws2.Range("A38").Value = "8th Attempt"
$ws2.Range("B38").Value = "Linear Regression"
$ws2.Range("C38").Value = 57319450.73
$ws2.Range("D38").Value = 81629775.77
$ws2.Range("E38").Value = 0.83
$ws2.Range("F38").Value = 0.75
$ws2.Range("G38").Value = "MLAssignment2_8"

# Row 39: Polynomial (no numeric results, left blank)
$ws2.Range("A39").Value = "8th Attempt"
$ws2.Range("B39").Value = "Polynomial"

# Row 40: Ridge
$ws2.Range("A40").Value = "8th Attempt"
$ws2.Range("B40").Value = "Ridge"
$ws2.Range("C40").Value = 58311575.32
$ws2.Range("D40").Value = 80243320.76
$ws2.Range("E40").Value = 0.83
$ws2.Range("F40").Value = 0.76

# Row 41: Ridge W Normalization
$ws2.Range("A41").Value = "8th Attempt"
$ws2.Range("B41").Value = "Ridge W Normalization"
$ws2.Range("C41").Value = 74602382.70
$ws2.Range("D41").Value = 87764855.88
$ws2.Range("E41").Value = 0.78
$ws2.Range("F41").Value = 0.73

# Restore the view: scroll/select so the frozen pane shows column K onward
# and the active cell sits on the newly entered block.
$ws2.Activate()
$ws2.Range("G42").Select()

# --- Sheet "Assign 1": move the selection ---
$ws1 = $wb.Worksheets.Item("Assign 1")
$ws1.Range("E64").Select()

# Leave "Assign 2" as the active sheet/tab, matching the saved workbook state.
$ws2.Activate()
